$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Bold "Please note..." paragraph: append " (Questions 1-5)" to
#    the existing run's text, then add a trailing run containing
#    just a space.
# ------------------------------------------------------------------
$noteRng = $d.Content
$noteRng.Find.Execute("Please note that the steps show rounded numbers, but that the final answers to the problems are calculated without rounding.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$noteRng.Collapse(0)
$noteRng.InsertAfter(" (Questions 1-5)")
$noteRng.Collapse(0)
$noteRng.InsertAfter(" ")

# ------------------------------------------------------------------
# 2. Remove the "Mode" row (Problem 1 / Mode / The most frequently
#    occurring value) from the solutions table entirely.
# ------------------------------------------------------------------
$t = $d.Tables.Item(1)
$t.Rows.Item(4).Delete()

# ------------------------------------------------------------------
# 3. Strip the trailing " Mode -XX.XXXXXX" text from each of the
#    Company A-E solution cells, leaving just Mean/Median.
# ------------------------------------------------------------------
$solutionTexts = @(
    "Mean: 21.276 Median: 13.433",
    "Mean: 33.482 Median: 20.838",
    "Mean: 41.122 Median: 25.558",
    "Mean: 0.706 Median: 1.892",
    "Mean: -1.084 Median: -3.796"
)
for ($i = 0; $i -lt $solutionTexts.Count; $i++) {
    $cell = $t.Rows.Item(4 + $i).Cells.Item(3)
    $cellRng = $cell.Range
    $cellRng.MoveEnd(1, -1)
    $cellRng.Text = $solutionTexts[$i]
}
